# Daten aktualisiert am 2023-06-27
# Append the (re-)published ticker list to the end of column A, directly
# below the existing data (which currently ends at row 5732).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$tickers = @(
    "AAF","ABDN","ABF","ANTO","AUTO","AV","BARC","BATS","BDEV","BEZ",
    "BF.B","BKG","BNZL","BRBY","BRK.B","BT-A","CCH","CRDA","DCC","DGE",
    "ENT","EXPN","FCIT","FRAS","GLEN","HLMA","HSBA","HSX","IMB","IMI",
    "INF","ITRK","JMAT","KGF","LGEN","LLOY","LSEG","MNDI","MNG","OCDO",
    "PHNX","PSON","REL","RMV","RR","RS1","SBRY","SDR","SGRO","SKG",
    "SMDS","SMT","SN","SPX","SSE","STAN","STJ","ULVR","UU","WEIR","WTB"
)

# Find the first empty row below the current data in column A.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$startRow = $lastRow + 1

for ($i = 0; $i -lt $tickers.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $tickers[$i]
}
